# TestData.xlsx - Login sheet update
# - Browser value switched from "firefox" to "chrome"
# - New column F added: "GoogleSearch" header / "GitHub" value
# - Selection moved to the new F2 cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Browser value (E2): firefox -> chrome
$ws.Range("E2").Value = "chrome"

# Add the new column F: header in row 1, value in row 2
$ws.Range("F1").Value = "GoogleSearch"
$ws.Range("F2").Value = "GitHub"

# Match F2's formatting reset (no visible border/fill, format explicitly applied)
$ws.Range("F2").Interior.Pattern = -4142
$ws.Range("F2").Borders.Item(7).LineStyle = -4142
$ws.Range("F2").Borders.Item(8).LineStyle = -4142
$ws.Range("F2").Borders.Item(9).LineStyle = -4142
$ws.Range("F2").Borders.Item(10).LineStyle = -4142

# Size the new column to fit its contents, like the other data columns
$ws.Columns.Item(6).EntireColumn.AutoFit()

# Move the active selection to the newly populated cell
$ws.Range("F2").Select() | Out-Null
